$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 77, pushing the existing rows 77..122 down to 78..123.
$ws.Rows(77).Insert()

# Populate the newly inserted row 77 with the new weekly record.
$ws.Cells.Item(77, 1).Value = 8
$ws.Cells.Item(77, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(77, 3).Value = "Coquimbo"
$ws.Cells.Item(77, 4).Value = 44609
$ws.Cells.Item(77, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(77, 5).Value = 4
$ws.Cells.Item(77, 6).Value = 100112044
$ws.Cells.Item(77, 7).Value = "Perejil"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 2600
$ws.Cells.Item(77, 11).Value = 2300
$ws.Cells.Item(77, 12).Value = 2500
$ws.Cells.Item(77, 13).Value = 2400
$ws.Cells.Item(77, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(77, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(77, 16).Value = 1600
$ws.Cells.Item(77, 17).Value = 1.5
$ws.Cells.Item(77, 18).Value = "Hortaliza"
